$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VillaTest")

# Replace the hotel name for Dubai (B3) with the new value.
$ws.Range("B3").Value = "Marco Polo Hotel"

# Move the active selection to B4, mirroring the post-edit cursor move.
$ws.Activate()
$ws.Range("B4").Select()
